# Update "想去人数" (column F) values across the sheets as per the
# upstream data refresh (commit: "Update gh-pages to output generated at 456a3b4").
#
# Sheet order in the workbook:
#   1 -> 展览     (Exhibitions)
#   2 -> 演出     (Performances)
#   3 -> 本地生活 (Local Life)
#   4 -> 全部类型 (All Types)

$wb = $excel.ActiveWorkbook

# ---- Sheet 1: 展览 ----
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value  = 14
$ws1.Range("F3").Value  = 346
$ws1.Range("F4").Value  = 1338
$ws1.Range("F5").Value  = 385
$ws1.Range("F6").Value  = 369
$ws1.Range("F7").Value  = 3939
$ws1.Range("F9").Value  = 791
$ws1.Range("F10").Value = 2386
$ws1.Range("F16").Value = 203
$ws1.Range("F17").Value = 3258
$ws1.Range("F19").Value = 238
$ws1.Range("F20").Value = 51
$ws1.Range("F21").Value = 353
$ws1.Range("F22").Value = 250
$ws1.Range("F23").Value = 56
$ws1.Range("F24").Value = 285

# ---- Sheet 2: 演出 ----
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F10").Value = 102
$ws2.Range("F11").Value = 230
$ws2.Range("F22").Value = 83

# ---- Sheet 3: 本地生活 ----
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F4").Value = 2130
$ws3.Range("F6").Value = 23

# ---- Sheet 4: 全部类型 ----
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F4").Value  = 2130
$ws4.Range("F8").Value  = 14
$ws4.Range("F10").Value = 346
$ws4.Range("F11").Value = 1338
$ws4.Range("F12").Value = 385
$ws4.Range("F16").Value = 23
$ws4.Range("F17").Value = 369
$ws4.Range("F18").Value = 3939
$ws4.Range("F23").Value = 102
$ws4.Range("F24").Value = 791
$ws4.Range("F25").Value = 2387
$ws4.Range("F28").Value = 230
$ws4.Range("F32").Value = 203
$ws4.Range("F38").Value = 238
$ws4.Range("F39").Value = 51
$ws4.Range("F40").Value = 353
$ws4.Range("F41").Value = 250
$ws4.Range("F42").Value = 56
$ws4.Range("F49").Value = 83
$ws4.Range("F50").Value = 285
